{"js": "// The document had a handful of leftover review artifacts from the lab's\n// editing pass: one tracked deletion (a stray \"-\" between \"Does \" and\n// \"diversity in\"), three tracked formatting changes (rPrChange records on\n// two \"sp.\" runs and one \"is\" run that had been italic and were changed to\n// non-italic), and five reviewer comments (with their anchoring\n// commentRangeStart/commentRangeEnd/commentReference markup) from\n// Mendelson, Alex and Morales, Breana M. The commit simply finalizes the\n// document for the undergrad's lab-meeting presentation: accept every\n// tracked change and remove every comment, leaving the plain, clean text\n// behind (no visible wording changes beyond what the tracked deletion\n// already encoded).\n\n// 1) Accept all tracked changes (the w:del around the stray hyphen and the\n//    w:rPrChange entries on the \"is\" / \"sp.\" runs) so the document reflects\n//    only the final, already-applied text/formatting.\nconst trackedChanges = context.document.body.getTrackedChanges();\ntrackedChanges.acceptAll();\nawait context.sync();\n\n// 2) Remove every comment (and its comment-range anchors) left in the\n//    document by the reviewers.\ncontext.document.deleteAllComments();\nawait context.sync();\n", "ps1": "# The document had a handful of leftover review artifacts from the lab's\n# editing pass: one tracked deletion (a stray \"-\" between \"Does \" and\n# \"diversity in\"), three tracked formatting changes (rPrChange records on\n# two \"sp.\" runs and one \"is\" run that had been italic and were changed to\n# non-italic), and five reviewer comments (with their anchoring\n# commentRangeStart/commentRangeEnd/commentReference markup) from\n# Mendelson, Alex and Morales, Breana M. The commit simply finalizes the\n# document for the undergrad's lab-meeting presentation: accept every\n# tracked change and remove every comment, leaving the plain, clean text\n# behind (no visible wording changes beyond what the tracked deletion\n# already encoded).\n\n$d = $word.ActiveDocument\n\n# 1) Accept all tracked changes (the deleted hyphen and the rPrChange\n#    formatting records on \"is\"/\"sp.\") so the document reflects only the\n#    final, already-applied text/formatting.\n$d.Revisions.AcceptAll()\n\n# 2) Remove every comment (and its comment-range anchors) left behind by\n#    the reviewers. Deleting from the front while the count shrinks avoids\n#    any index shifting issues.\nwhile ($d.Comments.Count -gt 0) {\n    $d.Comments.Item(1).Delete()\n}\n"}
